# Auto-update draw results: append the latest Pick 4 draw as a new row
# at the bottom of the results table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$newRow = $used.Row + $used.Rows.Count

# Leading apostrophes force Excel to store these values as plain text
# (matching the existing rows, which are all text cells) instead of
# silently converting them to dates/numbers.
$ws.Range("A$newRow").Value = "'2025-09-21"
$ws.Range("B$newRow").Value = "Pick 4"
$ws.Range("C$newRow").Value = "'250921"
$ws.Range("D$newRow").Value = "7-2-3-1"
$ws.Range("E$newRow").Value = "2025-09-21T21:35:40.293+04:00"
